# Update TPM-derived NATMI ligand-receptor metrics on Sheet1 (rows 2-9)
# to reflect the new TPM values, per "update scripts wuth new tpm".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 2.11624
$ws.Cells.Item(2, 8).Value = 6.34872
$ws.Cells.Item(2, 9).Value = 0.1897594766532197
$ws.Cells.Item(2, 10).Value = 0.1897594766532197
$ws.Cells.Item(2, 15).Value = 0.7091726973716084
$ws.Cells.Item(2, 16).Value = 0.7091726973716084
$ws.Cells.Item(2, 17).Value = 3.614086455466666
$ws.Cells.Item(2, 18).Value = 32.5267780992
$ws.Cells.Item(2, 19).Value = 0.1345722399099886
$ws.Cells.Item(2, 20).Value = 0.1345722399099886

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 2.11624
$ws.Cells.Item(3, 8).Value = 6.34872
$ws.Cells.Item(3, 9).Value = 0.1897594766532197
$ws.Cells.Item(3, 10).Value = 0.1897594766532197
$ws.Cells.Item(3, 13).Value = 0.7003526666666667
$ws.Cells.Item(3, 14).Value = 2.101058
$ws.Cells.Item(3, 15).Value = 0.2908273026283917
$ws.Cells.Item(3, 16).Value = 0.2908273026283917
$ws.Cells.Item(3, 17).Value = 1.482114327306667
$ws.Cells.Item(3, 18).Value = 13.33902894576
$ws.Cells.Item(3, 19).Value = 0.05518723674323116
$ws.Cells.Item(3, 20).Value = 0.05518723674323116

$ws.Cells.Item(4, 9).Value = 0.6160274054778138
$ws.Cells.Item(4, 10).Value = 0.6160274054778138
$ws.Cells.Item(4, 15).Value = 0.7091726973716084
$ws.Cells.Item(4, 16).Value = 0.7091726973716084
$ws.Cells.Item(4, 19).Value = 0.4368698167975347
$ws.Cells.Item(4, 20).Value = 0.4368698167975347

$ws.Cells.Item(5, 9).Value = 0.6160274054778138
$ws.Cells.Item(5, 10).Value = 0.6160274054778138
$ws.Cells.Item(5, 13).Value = 0.7003526666666667
$ws.Cells.Item(5, 14).Value = 2.101058
$ws.Cells.Item(5, 15).Value = 0.2908273026283917
$ws.Cells.Item(5, 16).Value = 0.2908273026283917
$ws.Cells.Item(5, 17).Value = 4.811475346450001
$ws.Cells.Item(5, 18).Value = 43.30327811805
$ws.Cells.Item(5, 19).Value = 0.1791575886802791
$ws.Cells.Item(5, 20).Value = 0.1791575886802791

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.793503666666667
$ws.Cells.Item(6, 8).Value = 5.380511
$ws.Cells.Item(6, 9).Value = 0.1608202836929164
$ws.Cells.Item(6, 10).Value = 0.1608202836929164
$ws.Cells.Item(6, 15).Value = 0.7091726973716084
$ws.Cells.Item(6, 16).Value = 0.7091726973716084
$ws.Cells.Item(6, 17).Value = 3.062921648551111
$ws.Cells.Item(6, 18).Value = 27.56629483696
$ws.Cells.Item(6, 19).Value = 0.1140493543785728
$ws.Cells.Item(6, 20).Value = 0.1140493543785728

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.793503666666667
$ws.Cells.Item(7, 8).Value = 5.380511
$ws.Cells.Item(7, 9).Value = 0.1608202836929164
$ws.Cells.Item(7, 10).Value = 0.1608202836929164
$ws.Cells.Item(7, 13).Value = 0.7003526666666667
$ws.Cells.Item(7, 14).Value = 2.101058
$ws.Cells.Item(7, 15).Value = 0.2908273026283917
$ws.Cells.Item(7, 16).Value = 0.2908273026283917
$ws.Cells.Item(7, 17).Value = 1.256085075626445
$ws.Cells.Item(7, 18).Value = 11.304765680638
$ws.Cells.Item(7, 19).Value = 0.0467709293143436
$ws.Cells.Item(7, 20).Value = 0.0467709293143436

$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.3724043333333333
$ws.Cells.Item(8, 8).Value = 1.117213
$ws.Cells.Item(8, 9).Value = 0.03339283417605023
$ws.Cells.Item(8, 10).Value = 0.03339283417605023
$ws.Cells.Item(8, 15).Value = 0.7091726973716084
$ws.Cells.Item(8, 16).Value = 0.7091726973716084
$ws.Cells.Item(8, 17).Value = 0.6359871550755555
$ws.Cells.Item(8, 18).Value = 5.72388439568
$ws.Cells.Item(8, 19).Value = 0.02368128628551237
$ws.Cells.Item(8, 20).Value = 0.02368128628551237

$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.3724043333333333
$ws.Cells.Item(9, 8).Value = 1.117213
$ws.Cells.Item(9, 9).Value = 0.03339283417605023
$ws.Cells.Item(9, 10).Value = 0.03339283417605023
$ws.Cells.Item(9, 13).Value = 0.7003526666666667
$ws.Cells.Item(9, 14).Value = 2.101058
$ws.Cells.Item(9, 15).Value = 0.2908273026283917
$ws.Cells.Item(9, 16).Value = 0.2908273026283917
$ws.Cells.Item(9, 17).Value = 0.2608143679282223
$ws.Cells.Item(9, 18).Value = 2.347329311354
$ws.Cells.Item(9, 19).Value = 0.00971154789053786
$ws.Cells.Item(9, 20).Value = 0.00971154789053786
